$d = $word.ActiveDocument

# 1. Title / heading text (appears twice - main Heading1 and bold summary near bottom)
$d.Content.Find.Execute(
    "Play Christmas Luck for Free - Unique One-Reel Slot with Exciting Features",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Christmas Luck Free: Exciting Features and Big Win Opportunities",
    2)

# 2. "What we like" bullet - unique gameplay
$d.Content.Find.Execute(
    "Unique gameplay with only one reel",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Unique one-reel feature adds excitement",
    2)

# 3. "What we like" bullet - bonus features
$d.Content.Find.Execute(
    "Exciting bonus features and symbols with great potential for big wins",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Special features and symbols offer big win opportunities",
    2)

# 4. "What we like" bullet - volatility
$d.Content.Find.Execute(
    "High volatility adds to the excitement of the game",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "High volatility for thrill-seeking players",
    2)

# 5. "What we like" bullet - demo version
$d.Content.Find.Execute(
    "Demo version available for free play to master the game",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Demo version available for free practice",
    2)

# 6. "What we don't like" bullet - RTP
$d.Content.Find.Execute(
    "RTP is not the highest compared to other slot games",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Lower RTP compared to other games",
    2)

# 7. "What we don't like" bullet - one reel
$d.Content.Find.Execute(
    "Only one reel may not appeal to players who prefer more reels and pay lines",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Limited reel may not appeal to traditional slot players",
    2)

# 8. Italic meta description near bottom
$d.Content.Find.Execute(
    "Read our review of Christmas Luck slot game and play for free. Unique one-reel gameplay with exciting features and high volatility for big wins.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Christmas Luck for free and enjoy its unique features and high volatility.",
    2)
